$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2986.7727
$ws.Range("I5").Value = 925.2941
$ws.Range("K5").Value = 925.2941
$ws.Range("M5").Value = -810.2941
$ws.Range("H34").Value = 1181.3334
$ws.Range("I34").Value = 1181.3334
$ws.Range("K34").Value = 1181.3334
$ws.Range("M34").Value = -978.3334
$ws.Range("H36").Value = 1181.3334
$ws.Range("I36").Value = 1181.3334
$ws.Range("K36").Value = 1181.3334
$ws.Range("M36").Value = -466.3334
$ws.Range("H86").Value = 100001270
$ws.Range("I86").Value = 100001270
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 100001270
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -100000147
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 100001270
$ws.Range("I89").Value = 100001270
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 500006350
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -500000734
$ws.Range("N89").ClearContents()
$ws.Range("H112").Value = 3219.2
$ws.Range("J112").Value = 1900.1305
$ws.Range("L112").Value = 5700.3915
$ws.Range("N112").Value = -7916.3915
$ws.Range("H135").Value = 3135.625
$ws.Range("I135").Value = 1861.25
$ws.Range("J135").Value = 4410
$ws.Range("K135").Value = 16751.25
$ws.Range("L135").Value = 39690
$ws.Range("M135").Value = -14216.25
$ws.Range("N135").Value = -44760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 481.9
$ws.Range("I4").Value = 461.375
$ws.Range("K4").Value = 461.375
$ws.Range("M4").Value = -345.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 7103.5
$ws.Range("I24").Value = 9499
$ws.Range("K24").Value = 9499
$ws.Range("M24").Value = -9264
$ws.Range("H64").Value = 4482.222
$ws.Range("I64").Value = 1547.75
$ws.Range("J64").Value = 6829.8
$ws.Range("K64").Value = 1547.75
$ws.Range("L64").Value = 6829.8
$ws.Range("M64").Value = -1322.75
$ws.Range("N64").Value = -7279.8
$ws.Range("H67").Value = 4482.222
$ws.Range("I67").Value = 1547.75
$ws.Range("J67").Value = 6829.8
$ws.Range("K67").Value = 1547.75
$ws.Range("L67").Value = 6829.8
$ws.Range("M67").Value = -767.75
$ws.Range("N67").Value = -8389.799999999999
$ws.Range("H80").Value = 15885838
$ws.Range("J80").Value = 22239276
$ws.Range("L80").Value = 22239276
$ws.Range("N80").Value = -22241272
$ws.Range("H83").Value = 15885838
$ws.Range("J83").Value = 22239276
$ws.Range("L83").Value = 111196380
$ws.Range("N83").Value = -111206364
$ws.Range("H94").Value = 2665.9443
$ws.Range("I94").Value = 1241.5
$ws.Range("K94").Value = 1241.5
$ws.Range("M94").Value = -790.5
$ws.Range("H107").Value = 16668416
$ws.Range("I107").Value = 20001496
$ws.Range("K107").Value = 20001496
$ws.Range("M107").Value = -19999576

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 360.92856
$ws.Range("I7").Value = 305.5
$ws.Range("J7").Value = 499.5
$ws.Range("K7").Value = 305.5
$ws.Range("L7").Value = 499.5
$ws.Range("M7").Value = -192.5
$ws.Range("N7").Value = -725.5
$ws.Range("H94").Value = 55559652
$ws.Range("I94").Value = 142859020
$ws.Range("J94").Value = 5507.5454
$ws.Range("K94").Value = 142859020
$ws.Range("L94").Value = 5507.5454
$ws.Range("M94").Value = -142858569
$ws.Range("N94").Value = -6409.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8539170
$ws.Range("I4").Value = 12500096
$ws.Range("K4").Value = 37500288
$ws.Range("M4").Value = -37500176
$ws.Range("H11").Value = 5312.7144
$ws.Range("I11").Value = 6039
$ws.Range("J11").Value = 3497
$ws.Range("K11").Value = 18117
$ws.Range("L11").Value = 10491
$ws.Range("M11").Value = -17977
$ws.Range("N11").Value = -10771
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H26").Value = 453.16666
$ws.Range("I26").Value = 476.27274
$ws.Range("K26").Value = 1428.81822
$ws.Range("M26").Value = -1140.81822
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H51").Value = 3999.6667
$ws.Range("I51").Value = 3999.6667
$ws.Range("K51").Value = 11999.0001
$ws.Range("M51").Value = -11539.0001
$ws.Range("H57").Value = 3331.6667
$ws.Range("I57").Value = 1999
$ws.Range("J57").Value = 9995
$ws.Range("K57").Value = 5997
$ws.Range("L57").Value = 29985
$ws.Range("M57").Value = -5438
$ws.Range("N57").Value = -31103
$ws.Range("H106").Value = 23400
$ws.Range("J106").Value = 23400
$ws.Range("L106").Value = 70200
$ws.Range("N106").Value = -72092
$ws.Range("H113").Value = 358052.16
$ws.Range("I113").Value = 497.5
$ws.Range("J113").Value = 501074
$ws.Range("K113").Value = 1492.5
$ws.Range("L113").Value = 1503222
$ws.Range("M113").Value = 677.5
$ws.Range("N113").Value = -1507562
$ws.Range("H134").Value = 5595.467
$ws.Range("I134").Value = 5595.467
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 16786.401
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11716.401
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 70
$ws.Range("K2").Value = 70
$ws.Range("M2").Value = 43
$ws.Range("H70").Value = 8360.423000000001
$ws.Range("I70").Value = 7852.3335
$ws.Range("J70").Value = 8795.929
$ws.Range("K70").Value = 7852.3335
$ws.Range("L70").Value = 8795.929
$ws.Range("M70").Value = -7582.3335
$ws.Range("N70").Value = -9335.929
$ws.Range("H73").Value = 8360.423000000001
$ws.Range("I73").Value = 7852.3335
$ws.Range("J73").Value = 8795.929
$ws.Range("K73").Value = 7852.3335
$ws.Range("L73").Value = 8795.929
$ws.Range("M73").Value = -6916.3335
$ws.Range("N73").Value = -10667.929
$ws.Range("H133").Value = 136048
$ws.Range("J133").Value = 136048
$ws.Range("L133").Value = 136048
$ws.Range("N133").Value = -146168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1042.4
$ws.Range("I32").Value = 1042.4
$ws.Range("K32").Value = 1042.4
$ws.Range("M32").Value = -725.4000000000001
$ws.Range("H46").Value = 25001910
$ws.Range("J46").Value = 50003000
$ws.Range("L46").Value = 50003000
$ws.Range("N46").Value = -50003376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4477
$ws.Range("I132").Value = 2682.2632
$ws.Range("K132").Value = 8046.7896
$ws.Range("M132").Value = -5516.7896
$ws.Range("H139").Value = 84968.5
$ws.Range("J139").Value = 84968.5
$ws.Range("L139").Value = 84968.5
$ws.Range("N139").Value = -95248.5
